$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Permits Filed for 8802 162nd Street in Jamaica, Queens"
$ws.Range("B2").Value = "https://newyorkyimby.com/2025/12/8802-162nd-street-in-jamaica-queens.html"
$ws.Range("C2").Value = 'Permits have been filed for a 14-story mixed-use building at 8802 162nd Street in <a href="https://newyorkyimby.com/neighborhoods/jamaica">Jamaica</a>, Queens. Located at the intersection of Highland Avenue and 162nd Street, the lot is near the Parsons Boulevard subway station, served by the E and F trains. Sandi Silk of Mural Real Estate Partners is listed as the owner behind the applications.'
$ws.Range("D2").Value = "2025-12-16T11:30:21+00:00"
$ws.Range("E2").Value = "Tue, 16 Dec 2025 11:30:21 +0000"
$ws.Range("F2").Value = "YIMBY"
$ws.Range("G2").Value = "YIMBY - Jamaica"
$ws.Range("H2").Value = ""
